# Auto-generated edit script applying the Diabolos_Profits.xlsx market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Cells.Item(19,8).Value = 14131.667
$ws.Cells.Item(19,9).Value = 1887.7
$ws.Cells.Item(19,10).Value = 25262.545
$ws.Cells.Item(19,11).Value = 1887.7
$ws.Cells.Item(19,12).Value = 25262.545
$ws.Cells.Item(19,13).Value = -1712.7
$ws.Cells.Item(19,14).Value = -25612.545
# Row 28
$ws.Cells.Item(28,8).Value = 26810.105
$ws.Cells.Item(28,9).Value = 32789.195
$ws.Cells.Item(28,11).Value = 32789.195
$ws.Cells.Item(28,13).Value = -32304.195
# Row 40
$ws.Cells.Item(40,8).Value = 3100
$ws.Cells.Item(40,10).Value = 3100
$ws.Cells.Item(40,12).Value = 3100
$ws.Cells.Item(40,14).Value = -3450
# Row 137
$ws.Cells.Item(137,8).Value = 2481.8293
$ws.Cells.Item(137,9).Value = 2267.0967
$ws.Cells.Item(137,10).Value = 3147.5
$ws.Cells.Item(137,11).Value = 6801.2901
$ws.Cells.Item(137,12).Value = 9442.5
$ws.Cells.Item(137,13).Value = -4251.2901
$ws.Cells.Item(137,14).Value = -14542.5
# Row 138
$ws.Cells.Item(138,8).Value = 1819.5916
$ws.Cells.Item(138,9).Value = 1437.2069
$ws.Cells.Item(138,10).Value = 2083.6191
$ws.Cells.Item(138,11).Value = 4311.620699999999
$ws.Cells.Item(138,12).Value = 6250.8573
$ws.Cells.Item(138,13).Value = 828.3793000000005
$ws.Cells.Item(138,14).Value = -16530.8573
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2,8).Value = 828.2222
$ws.Cells.Item(2,9).Value = 882.0625
$ws.Cells.Item(2,10).Value = 397.5
$ws.Cells.Item(2,11).Value = 882.0625
$ws.Cells.Item(2,12).Value = 397.5
$ws.Cells.Item(2,13).Value = -769.0625
$ws.Cells.Item(2,14).Value = -623.5
# Row 45
$ws.Cells.Item(45,8).Value = 2118.75
$ws.Cells.Item(45,9).Value = 1700
$ws.Cells.Item(45,11).Value = 1700
$ws.Cells.Item(45,13).Value = -1323
# Row 61
$ws.Cells.Item(61,8).Value = 2719.6667
$ws.Cells.Item(61,9).Value = 1813.7
$ws.Cells.Item(61,11).Value = 1813.7
$ws.Cells.Item(61,13).Value = -1601.7
# Row 63
$ws.Cells.Item(63,8).Value = 60004164
$ws.Cells.Item(63,9).Value = 83335310
$ws.Cells.Item(63,11).Value = 83335310
$ws.Cells.Item(63,13).Value = -83334624
# Row 66
$ws.Cells.Item(66,8).Value = 60004164
$ws.Cells.Item(66,9).Value = 83335310
$ws.Cells.Item(66,11).Value = 416676550
$ws.Cells.Item(66,13).Value = -416673118
# Row 88
$ws.Cells.Item(88,8).Value = 18520014
$ws.Cells.Item(88,10).Value = 1959.6
$ws.Cells.Item(88,12).Value = 1959.6
$ws.Cells.Item(88,14).Value = -2771.6
# Row 91
$ws.Cells.Item(91,8).Value = 18520014
$ws.Cells.Item(91,10).Value = 1959.6
$ws.Cells.Item(91,12).Value = 1959.6
$ws.Cells.Item(91,14).Value = -4767.6
# Row 116
$ws.Cells.Item(116,8).Value = 828.2222
$ws.Cells.Item(116,9).Value = 882.0625
$ws.Cells.Item(116,10).Value = 397.5
$ws.Cells.Item(116,11).Value = 882.0625
$ws.Cells.Item(116,12).Value = 397.5
$ws.Cells.Item(116,13).Value = 1411.9375
$ws.Cells.Item(116,14).Value = -4985.5
# Row 136
$ws.Cells.Item(136,8).Value = 2719.6667
$ws.Cells.Item(136,9).Value = 1813.7
$ws.Cells.Item(136,11).Value = 5441.1
$ws.Cells.Item(136,13).Value = -2891.1
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3,8).Value = 828.2222
$ws.Cells.Item(3,9).Value = 882.0625
$ws.Cells.Item(3,10).Value = 397.5
$ws.Cells.Item(3,11).Value = 882.0625
$ws.Cells.Item(3,12).Value = 397.5
$ws.Cells.Item(3,13).Value = -768.0625
$ws.Cells.Item(3,14).Value = -625.5
# Row 20
$ws.Cells.Item(20,8).Value = 17581.875
$ws.Cells.Item(20,9).Value = 24572.908
$ws.Cells.Item(20,10).Value = 2201.6
$ws.Cells.Item(20,11).Value = 24572.908
$ws.Cells.Item(20,12).Value = 2201.6
$ws.Cells.Item(20,13).Value = -24325.908
$ws.Cells.Item(20,14).Value = -2695.6
# Row 22
$ws.Cells.Item(22,8).Value = 312.16666
$ws.Cells.Item(22,9).Value = 312.16666
$ws.Cells.Item(22,11).Value = 312.16666
$ws.Cells.Item(22,13).Value = -139.16666
# Row 105
$ws.Cells.Item(105,8).Value = 2036.1714
$ws.Cells.Item(105,9).Value = 1792.875
$ws.Cells.Item(105,10).Value = 2567
$ws.Cells.Item(105,11).Value = 1792.875
$ws.Cells.Item(105,12).Value = 2567
$ws.Cells.Item(105,13).Value = -45.875
$ws.Cells.Item(105,14).Value = -6061
# Row 107
$ws.Cells.Item(107,8).Value = 50201908
$ws.Cells.Item(107,9).Value = 335262.66
$ws.Cells.Item(107,10).Value = 125001870
$ws.Cells.Item(107,11).Value = 335262.66
$ws.Cells.Item(107,12).Value = 125001870
$ws.Cells.Item(107,13).Value = -333342.66
$ws.Cells.Item(107,14).Value = -125005710
# Row 132
$ws.Cells.Item(132,8).Value = 0
$ws.Cells.Item(132,10).Value = 0
$ws.Cells.Item(132,12).Value = 0
$ws.Cells.Item(132,14).ClearContents()
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16,8).Value = 1441.8667
$ws.Cells.Item(16,9).Value = 1259.1428
$ws.Cells.Item(16,11).Value = 1259.1428
$ws.Cells.Item(16,13).Value = -972.1428000000001
# Row 50
$ws.Cells.Item(50,8).Value = 37500
$ws.Cells.Item(50,9).Value = 37500
$ws.Cells.Item(50,11).Value = 37500
$ws.Cells.Item(50,13).Value = -36875
# Row 76
$ws.Cells.Item(76,8).Value = 14999
$ws.Cells.Item(76,9).Value = 14999
$ws.Cells.Item(76,11).Value = 14999
$ws.Cells.Item(76,13).Value = -14684
# Row 79
$ws.Cells.Item(79,8).Value = 14999
$ws.Cells.Item(79,9).Value = 14999
$ws.Cells.Item(79,11).Value = 14999
$ws.Cells.Item(79,13).Value = -13907
# Row 107
$ws.Cells.Item(107,8).Value = 742.26666
$ws.Cells.Item(107,9).Value = 722.9286
$ws.Cells.Item(107,10).Value = 1013
$ws.Cells.Item(107,11).Value = 722.9286
$ws.Cells.Item(107,12).Value = 1013
$ws.Cells.Item(107,13).Value = 1197.0714
$ws.Cells.Item(107,14).Value = -4853
# Row 113
$ws.Cells.Item(113,8).Value = 1441.8667
$ws.Cells.Item(113,9).Value = 1259.1428
$ws.Cells.Item(113,11).Value = 1259.1428
$ws.Cells.Item(113,13).Value = 910.8571999999999
# Row 132
$ws.Cells.Item(132,8).Value = 1802.8572
$ws.Cells.Item(132,9).Value = 1802.8572
$ws.Cells.Item(132,10).Value = 0
$ws.Cells.Item(132,11).Value = 5408.571599999999
$ws.Cells.Item(132,12).Value = 0
$ws.Cells.Item(132,13).Value = -2878.571599999999
$ws.Cells.Item(132,14).ClearContents()
$ws = $wb.Worksheets.Item("CUL")
# Row 40
$ws.Cells.Item(40,8).Value = 59.076923
$ws.Cells.Item(40,10).Value = 14
$ws.Cells.Item(40,12).Value = 56
$ws.Cells.Item(40,14).Value = -194
# Row 87
$ws.Cells.Item(87,8).Value = 21
$ws.Cells.Item(87,9).Value = 21
$ws.Cells.Item(87,11).Value = 63
$ws.Cells.Item(87,13).Value = 1185
# Row 90
$ws.Cells.Item(90,8).Value = 21
$ws.Cells.Item(90,9).Value = 21
$ws.Cells.Item(90,11).Value = 189
$ws.Cells.Item(90,13).Value = 6051
# Row 103
$ws.Cells.Item(103,8).Value = 869.625
$ws.Cells.Item(103,9).Value = 674.25
$ws.Cells.Item(103,10).Value = 1065
$ws.Cells.Item(103,11).Value = 2022.75
$ws.Cells.Item(103,12).Value = 3195
$ws.Cells.Item(103,13).Value = -1143.75
$ws.Cells.Item(103,14).Value = -4953
# Row 107
$ws.Cells.Item(107,8).Value = 443.36365
$ws.Cells.Item(107,9).Value = 614.5
$ws.Cells.Item(107,10).Value = 238
$ws.Cells.Item(107,11).Value = 1843.5
$ws.Cells.Item(107,12).Value = 714
$ws.Cells.Item(107,13).Value = 76.5
$ws.Cells.Item(107,14).Value = -4554
# Row 114
$ws.Cells.Item(114,8).Value = 3531.5625
$ws.Cells.Item(114,10).Value = 5939
$ws.Cells.Item(114,12).Value = 17817
$ws.Cells.Item(114,14).Value = -24325
# Row 139
$ws.Cells.Item(139,8).Value = 25001742
$ws.Cells.Item(139,9).Value = 33334902
$ws.Cells.Item(139,10).Value = 2257.4
$ws.Cells.Item(139,11).Value = 100004706
$ws.Cells.Item(139,12).Value = 6772.200000000001
$ws.Cells.Item(139,13).Value = -99999566
$ws.Cells.Item(139,14).Value = -17052.2
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Cells.Item(70,8).Value = 7205.316
$ws.Cells.Item(70,9).Value = 6300.846
$ws.Cells.Item(70,10).Value = 9165
$ws.Cells.Item(70,11).Value = 6300.846
$ws.Cells.Item(70,12).Value = 9165
$ws.Cells.Item(70,13).Value = -6030.846
$ws.Cells.Item(70,14).Value = -9705
# Row 73
$ws.Cells.Item(73,8).Value = 7205.316
$ws.Cells.Item(73,9).Value = 6300.846
$ws.Cells.Item(73,10).Value = 9165
$ws.Cells.Item(73,11).Value = 6300.846
$ws.Cells.Item(73,12).Value = 9165
$ws.Cells.Item(73,13).Value = -5364.846
$ws.Cells.Item(73,14).Value = -11037
# Row 113
$ws.Cells.Item(113,8).Value = 2420
$ws.Cells.Item(113,9).Value = 1760.8695
$ws.Cells.Item(113,11).Value = 1760.8695
$ws.Cells.Item(113,13).Value = 409.1305
# Row 136
$ws.Cells.Item(136,8).Value = 22378.45
$ws.Cells.Item(136,10).Value = 22378.45
$ws.Cells.Item(136,12).Value = 67135.35000000001
$ws.Cells.Item(136,14).Value = -72235.35000000001
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Cells.Item(132,8).Value = 43481240
$ws.Cells.Item(132,9).Value = 55558228
$ws.Cells.Item(132,10).Value = 4092.8
$ws.Cells.Item(132,11).Value = 166674684
$ws.Cells.Item(132,12).Value = 12278.4
$ws.Cells.Item(132,13).Value = -166672154
$ws.Cells.Item(132,14).Value = -17338.4
# Row 136
$ws.Cells.Item(136,8).Value = 3174.3076
$ws.Cells.Item(136,10).Value = 5484.5
$ws.Cells.Item(136,12).Value = 16453.5
$ws.Cells.Item(136,14).Value = -21553.5
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Cells.Item(81,8).Value = 22227608
$ws.Cells.Item(81,9).Value = 5620
$ws.Cells.Item(81,10).Value = 40005200
$ws.Cells.Item(81,11).Value = 11240
$ws.Cells.Item(81,12).Value = 80010400
$ws.Cells.Item(81,13).Value = -10179
$ws.Cells.Item(81,14).Value = -80012522
# Row 84
$ws.Cells.Item(84,8).Value = 22227608
$ws.Cells.Item(84,9).Value = 5620
$ws.Cells.Item(84,10).Value = 40005200
$ws.Cells.Item(84,11).Value = 56200
$ws.Cells.Item(84,12).Value = 400052000
$ws.Cells.Item(84,13).Value = -50896
$ws.Cells.Item(84,14).Value = -400062608
# Row 107
$ws.Cells.Item(107,8).Value = 661.1515000000001
$ws.Cells.Item(107,9).Value = 647.5238000000001
$ws.Cells.Item(107,10).Value = 685
$ws.Cells.Item(107,11).Value = 1942.5714
$ws.Cells.Item(107,12).Value = 2055
$ws.Cells.Item(107,13).Value = -22.57140000000027
$ws.Cells.Item(107,14).Value = -5895
# Row 113
$ws.Cells.Item(113,8).Value = 6026.2607
$ws.Cells.Item(113,9).Value = 7670.7144
$ws.Cells.Item(113,10).Value = 3468.2222
$ws.Cells.Item(113,11).Value = 23012.1432
$ws.Cells.Item(113,12).Value = 10404.6666
$ws.Cells.Item(113,13).Value = -20842.1432
$ws.Cells.Item(113,14).Value = -14744.6666
# Row 122
$ws.Cells.Item(122,8).Value = 2242.4092
$ws.Cells.Item(122,10).Value = 2275.5
$ws.Cells.Item(122,12).Value = 6826.5
$ws.Cells.Item(122,14).Value = -11726.5
# Row 132
$ws.Cells.Item(132,8).Value = 5226.8335
$ws.Cells.Item(132,9).Value = 4977.0625
$ws.Cells.Item(132,10).Value = 7225
$ws.Cells.Item(132,11).Value = 14931.1875
$ws.Cells.Item(132,12).Value = 21675
$ws.Cells.Item(132,13).Value = -12401.1875
$ws.Cells.Item(132,14).Value = -26735
